$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (kept as literal text, matching
# the source data which stores Price/Volume columns as text strings).
$updates = @{
    "D2" = "313.71"
    "E2" = "0.44%"
    "D3" = "48.67"
    "E3" = "9.90%"
    "D4" = "5.277"
    "E4" = "3.20%"
    "D5" = "0.07886"
    "D6" = "4.581"
    "E6" = "1.89%"
    "D7" = "1.318"
    "D8" = "1.612"
    "E8" = "-2.59%"
    "D9" = "0.1242"
    "E9" = "-4.21%"
    "D10" = "0.1931"
    "E10" = "1.28%"
    "D11" = "0.09417"
    "E11" = "0.36%"
    "D12" = "0.04545"
    "E12" = "6.98%"
    "D13" = "0.1047"
    "E13" = "0.73%"
    "D14" = "0.001305"
    "E14" = "-0.66%"
    "D15" = "0.04212"
    "E15" = "0.42%"
    "D16" = "0.005870"
    "E16" = "1.89%"
    "D17" = "3.337"
    "E17" = "-1.12%"
    "D18" = "2.467"
    "E18" = "2.70%"
    "E19" = "1.70%"
    "D20" = "8.051"
    "E20" = "-0.20%"
    "D21" = "0.1395"
    "E21" = "3.38%"
    "D22" = "0.3070"
    "E22" = "-1.89%"
    "D23" = "0.001296"
    "E23" = "1.84%"
    "E24" = "-8.87%"
    "D25" = "0.0001355"
    "E25" = "1.25%"
    "D26" = "0.0003556"
    "D38" = "0.02635"
    "E38" = "-0.78%"
    "D39" = "0.05809"
    "E39" = "7.26%"
    "D40" = "0.01080"
    "E40" = "92.15%"
    "D41" = "0.008009"
    "E41" = "3.67%"
    "D42" = "0.1443"
    "E42" = "2.28%"
    "D43" = "0.008457"
    "E43" = "15.41%"
    "D44" = "0.008666"
    "E44" = "10.13%"
    "D45" = "0.3132"
    "E45" = "0.29%"
    "D46" = "0.00006995"
    "E46" = "3.72%"
    "D47" = "0.00000000753"
    "E47" = "1.22%"
    "D48" = "0.05599"
    "E48" = "-0.02%"
    "D49" = "0.004019"
    "E49" = "1.35%"
    "D50" = "0.00002108"
    "E50" = "1.22%"
    "D51" = "0.0002007"
    "E51" = "1.22%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text interpretation so values like "313.71" or "0.44%" are not
    # auto-converted into numbers/percentages by Excel's smart-entry logic.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
